# Apply crypto price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.995.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.03%  '
$ws.Range("D3").Value = "'3.379.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'570.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").Value = "'140.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("E9").Value = '  +1.57%  '
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("E11").Value = '  -0.15%  '
$ws.Range("D12").Value = "'3.957.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("E13").Value = '  +2.08%  '
$ws.Range("D14").Value = "'27.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.55%  '
$ws.Range("D15").Value = "'3.377.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.06%  '
$ws.Range("E16").Value = '  -0.20%  '
$ws.Range("D17").Value = "'61.095.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E18").Value = '  -1.45%  '
$ws.Range("E19").Value = '  -1.81%  '
$ws.Range("D20").Value = "'8.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.34%  '
$ws.Range("D21").Value = "'381.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.02%  '
$ws.Range("D22").Value = "'75.33"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.22%  '
$ws.Range("E23").Value = '  -0.82%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").Value = "'3.515.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("E26").Value = '  -1.81%  '
$ws.Range("E27").Value = '  +8.32%  '
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("D29").Value = "'7.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.32%  '
$ws.Range("E30").Value = '  +0.51%  '
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("E33").Value = '  -3.68%  '
$ws.Range("D34").Value = "'23.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.06%  '
$ws.Range("D35").Value = "'6.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.77%  '
$ws.Range("D36").Value = "'166.04"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'3.414.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.26%  '
$ws.Range("E38").Value = '  +0.66%  '
$ws.Range("E39").Value = '  -1.94%  '
$ws.Range("E40").Value = '  -0.57%  '
$ws.Range("D41").Value = "'26.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.28%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("E44").Value = '  -1.46%  '
$ws.Range("E45").Value = '  -2.84%  '
$ws.Range("E46").Value = '  -0.52%  '
$ws.Range("D47").Value = "'2.434.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.21%  '
$ws.Range("D48").Value = "'22.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.37%  '
$ws.Range("E49").Value = '  -1.74%  '
$ws.Range("E50").Value = '  -2.35%  '
$ws.Range("D51").Value = "'2.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.72%  '
